# SmeltzPFAS-Clint-Level4.xlsx -- "updated files for marci's paper"
#
# Re-measured Clint values came back slightly different, the
# Clint.pValue computed for propranolol underflowed to 0 (so its
# scientific-notation number format is no longer needed), the
# Sat.pValue column was recomputed, and the table was re-styled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated Clint measurements (row 2 = Phenacetin) ---
$ws.Range("D2").Value2 = 10.7
$ws.Range("F2").Value2 = 14.2
$ws.Range("J2").Value2 = 0.0328

# --- Updated Clint measurements (row 3 = propranolol) ---
$ws.Range("D3").Value2 = 10
$ws.Range("E3").Value2 = 7.6
$ws.Range("F3").Value2 = 12.1
# Clint.pValue underflowed to 0 -- drop the custom 0.00E+00 number
# format it used to need and go back to the sheet's default style.
$ws.Range("J3").Style = "Normal"
$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 0.748

# --- Updated Sat.pValue column for the remaining rows ---
$ws.Range("K4").Value2 = 0.749

$ws.Range("J5").Value2 = 0.999
$ws.Range("K5").Value2 = 0.756

$ws.Range("K6").Value2 = 0.748

$ws.Range("K7").Value2 = 0.749

$ws.Range("K8").Value2 = 0.747

# --- Re-style the data table ---
$ws.ListObjects.Item(1).TableStyle = "TableStyleLight9"
